$wb = $excel.ActiveWorkbook

# Both "展览" (Exhibition) and "全部类型" (All Types) sheets list the same
# events and need the same "想去人数" (want-to-go count) updates in column F.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1218
    $ws.Range("F3").Value = 431
    $ws.Range("F4").Value = 298

    if ($sheetName -eq "展览") {
        $ws.Range("F7").Value = 12440
        $ws.Range("F13").Value = 12289
        $ws.Range("F14").Value = 4865
        $ws.Range("F15").Value = 4757
        $ws.Range("F21").Value = 7
    }
    else {
        $ws.Range("F9").Value = 12440
        $ws.Range("F15").Value = 12289
        $ws.Range("F16").Value = 4865
        $ws.Range("F17").Value = 4757
        $ws.Range("F23").Value = 7
    }
}
